$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the first small "Puntos de funcion" table (rows 15-19) ---
$ws.Range("C15").Value = 28.99
$ws.Range("C16").Value = 46.83
$ws.Range("C17").Value = 84.74
$ws.Range("C18").Value = 75.82
$ws.Range("C19").Value = 78.05

# --- Update the second table (rows 24-27) that feeds E/G formulas ---
$ws.Range("C24").Value = 46.83
$ws.Range("C25").Value = 84.74
$ws.Range("C26").Value = 75.82

# C27 used to hold the formula =19+13; it becomes a hard-coded value.
$ws.Range("C27").Value = 107.04

# --- Re-position the pie chart (was anchored near row 127, col H) ---
$chartObj = $ws.ChartObjects(1)
$chartObj.Left = 876.045
$chartObj.Top = 42
